# Updating scores for the 13th may
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("may18")

# Remove the stray "Check this" comment in F19 (shared string no longer used)
$ws.Range("F19").ClearContents()

# New round header (13 May 2018) starting at row 22
$ws.Range("A22").Value = (Get-Date -Year 2018 -Month 5 -Day 13 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("A22").NumberFormat = $ws.Range("A1").NumberFormat
$ws.Range("B22").Value = "Score"
$ws.Range("C22").Value = "Fairway"
$ws.Range("D22").Value = "GIR"
$ws.Range("E22").Value = "Putts"
$ws.Range("F22").Value = "Comment"

# Hole-by-hole data for the new round (rows 23-40)
$holes = @(
    @{ Hole = "Hole 1";  Score = 6; Fairway = "S"; Putts = 3 },
    @{ Hole = "Hole 2";  Score = 5; Fairway = "";  Putts = 2 },
    @{ Hole = "Hole 3";  Score = 4; Fairway = "S"; Putts = 2 },
    @{ Hole = "Hole 4";  Score = 5; Fairway = "S"; Putts = 2 },
    @{ Hole = "Hole 5";  Score = 4; Fairway = "";  Putts = 2 },
    @{ Hole = "Hole 6";  Score = 6; Fairway = "R"; Putts = 2 },
    @{ Hole = "Hole 7";  Score = 4; Fairway = "R"; Putts = 2 },
    @{ Hole = "Hole 8";  Score = 4; Fairway = "";  Putts = 2 },
    @{ Hole = "Hole 9";  Score = 5; Fairway = "S"; Putts = 3 },
    @{ Hole = "Hole 10"; Score = 4; Fairway = "R"; Putts = 2 },
    @{ Hole = "Hole 11"; Score = 3; Fairway = "";  Putts = 1 },
    @{ Hole = "Hole 12"; Score = 5; Fairway = "S"; Putts = 2 },
    @{ Hole = "Hole 13"; Score = 4; Fairway = "S"; Putts = 1 },
    @{ Hole = "Hole 14"; Score = 5; Fairway = "S"; Putts = 2 },
    @{ Hole = "Hole 15"; Score = 5; Fairway = "";  Putts = 2 },
    @{ Hole = "Hole 16"; Score = 4; Fairway = "L"; Putts = 1 },
    @{ Hole = "Hole 17"; Score = 7; Fairway = "R"; Putts = 2 },
    @{ Hole = "Hole 18"; Score = 4; Fairway = "S"; Putts = 2 }
)

$row = 23
foreach ($h in $holes) {
    $ws.Cells.Item($row, 1).Value = $h.Hole
    $ws.Cells.Item($row, 2).Value = $h.Score
    if ($h.Fairway -ne "") {
        $ws.Cells.Item($row, 3).Value = $h.Fairway
    }
    $ws.Cells.Item($row, 5).Value = $h.Putts
    $row++
}

# Totals row (row 41)
$ws.Range("B41").Formula = "=SUM(B23:B40)"
$ws.Range("E41").Formula = "=SUM(E23:E40)"

$ws.Range("C41").Select()
